$d = $word.ActiveDocument

# 1. Replace the single run "Основные методы:" with the new, longer text
#    "Другие основные основные методы:" (still a single run/formatting run
#    at this point - it gets split below once the bookmark is inserted).
$r = $d.Content
$null = $r.Find.Execute("Основные методы:", $false, $false, $false, $false, $false, $true, 1, $false, "Другие основные основные методы:", 2)

# 2. Find the split point inside the new text - right after
#    "Другие основные о" - which is where the _GoBack bookmark must sit,
#    splitting the run in two.
$r2 = $d.Content
$null = $r2.Find.Execute("Другие основные о", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
$splitPoint = $r2.End

# 3. Re-anchor the (hidden) _GoBack bookmark at that collapsed range. Adding
#    a bookmark with a name that already exists moves it there and removes
#    it from its previous location (after "x in s - принадлежит ли x
#    множеству s"), matching both halves of the diff in one step.
$splitRange = $d.Range($splitPoint, $splitPoint)
$d.Bookmarks.Add("_GoBack", $splitRange)
